$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.058.51"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.470.30"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.55"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.72"

$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.470.06"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("E10").Value = "  +0.45%  "

$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.72"
$ws.Range("E14").Value = "  +4.87%  "

$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.918.53"
$ws.Range("E16").Value = "  +2.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.062.60"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.472.29"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.20"
$ws.Range("E19").Value = "  +3.76%  "

$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.34"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.23"
$ws.Range("E22").Value = "  +9.40%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.27"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "668.43"
$ws.Range("E26").Value = "  +4.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  +12.92%  "

$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +438.30%  "

$ws.Range("E31").Value = "  +2.66%  "

$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("E35").Value = "  +3.41%  "

$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.47"
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.76"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.43"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0306"
$ws.Range("E45").Value = "  +4.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.83"
$ws.Range("E46").Value = "  +4.77%  "

$ws.Range("E47").Value = "  +19.30%  "

$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.62"
$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("E51").Value = "  -0.86%  "
